$wb = $excel.ActiveWorkbook

# --- Yearly sheet ---
$yearly = $wb.Worksheets.Item("Yearly")

# Update the raw dividend figures for February (row 4):
#   L4 = Taxable Account, N4 = Suzie's Roth IRA
$yearly.Range("L4").Value = 67.57
$yearly.Range("N4").Value = 22.05

# --- All Time sheet ---
$allTime = $wb.Worksheets.Item("All Time")

# Row 8 (year 2017) should pull its Taxable/401K/Suzie figures from the
# Yearly sheet's totals row (15) instead of a single month / hard values.
$allTime.Range("F8").Formula = "=Yearly!L15"
$allTime.Range("G8").Formula = "=Yearly!M15"
$allTime.Range("H8").Formula = "=Yearly!N15"

# --- Recalculate so cached formula values are refreshed ---
$excel.Calculate() | Out-Null

# --- Update the saved selections on each sheet ---
$yearly.Activate() | Out-Null
$yearly.Range("N15").Select() | Out-Null

$allTime.Activate() | Out-Null
$allTime.Range("L12").Select() | Out-Null
